$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K -> E:L
$ws.Columns("D").Insert()

# Copy cell formatting (number format/style) from the new right neighbor (E) into D
# so the new column inherits the same look (date format row 7/38/80, number format elsewhere).
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the newly reported period's figures
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 3096200
$ws.Range("D9").Value2 = 1695600
$ws.Range("D10").Value2 = 1400600
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = 278600
$ws.Range("D17").Value2 = 2460300
$ws.Range("D18").Value2 = 635800
$ws.Range("D20").Value2 = 22700
$ws.Range("D21").Value2 = 937100
$ws.Range("D22").Value2 = 151700
$ws.Range("D23").Value2 = 506800
$ws.Range("D24").Value2 = 129100
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 377700
$ws.Range("D27").Value2 = 372200
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -22700
$ws.Range("D33").Value2 = 372200
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 372200
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 87700
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 420100
$ws.Range("D44").Value2 = 44500
$ws.Range("D45").Value2 = 242400
$ws.Range("D46").Value2 = 794800
$ws.Range("D47").Value2 = 0
$ws.Range("D48").Value2 = 8909700
$ws.Range("D49").Value2 = 0
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 1705300
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 11409700
$ws.Range("D57").Value2 = 351400
$ws.Range("D58").Value2 = 250000
$ws.Range("D59").Value2 = 382000
$ws.Range("D60").Value2 = 983400
$ws.Range("D61").Value2 = 3120200
$ws.Range("D62").Value2 = 3106700
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 7210400
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 1727900
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 4199300
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 372200
$ws.Range("D83").Value2 = 278600
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 588100
$ws.Range("D91").Value2 = -864100
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -863900
$ws.Range("D96").Value2 = -65600
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = 357900
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = 82000
